# Add 2022-Q3 data:
#  1. Insert a new top data row in the "总计" (summary) sheet for 2022-Q3.
#  2. Insert a brand-new "2022-Q3" worksheet (positioned right after "总计")
#     holding the per-fund holdings detail for that quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" summary sheet - insert new row 2 with the 2022-Q3 totals.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 6
$summary.Cells.Item(2, 4).Value = 3.1

# The sheet's leftmost column is a re-computed 0-based row index (not a
# static label), so every row below the inserted one shifts up by one.
for ($r = 3; $r -le 9; $r++) {
    $summary.Cells.Item($r, 1).Value = $r - 2
}

# Row-Insert() only carried the bordered/centered "index-column" style onto
# B:D (copied down from the header row); re-stamp A2 with the same style
# used by the other index cells (A3, etc.) and drop the stray formatting
# that landed on B2:D2 so the new row matches its neighbours.
$summary.Range("B2:D2").ClearFormats()
$summary.Cells.Item(3, 1).Copy()
$summary.Cells.Item(2, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) Brand-new "2022-Q3" sheet, inserted right after "总计".
# ---------------------------------------------------------------------
$anchor = $wb.Worksheets.Item(2)
$q3 = $wb.Worksheets.Add($anchor)
$q3.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $q3.Cells.Item(1, $c + 2).Value = $headers[$c]
}

# Columns B-G hold text (fund code keeps leading zeros, numeric-looking
# figures keep their original trailing-zero formatting); only H (rank) is
# a real number. Force B:G to text before writing so Excel doesn't coerce
# values like "014834" or "8.50" into numbers.
$q3.Range("B2:G7").NumberFormat = "@"

$rows = @(
    @("014834", "汇添富盈鑫灵活配置混合D", "45.43", "91.88", "2.95", "1.3402", 10),
    @("002420", "汇添富盈鑫灵活配置混合A", "22.62", "91.88", "2.95", "0.6673", 10),
    @("257010", "国联安小盘精选混合",       "8.50",  "74.70", "6.98", "0.5933", 1),
    @("014833", "汇添富盈鑫灵活配置混合C", "14.93", "91.88", "2.95", "0.4404", 10),
    @("006138", "国联安价值优选股票",       "0.57",  "94.64", "6.82", "0.0389", 1),
    @("002367", "国联安安稳灵活配置混合",   "0.57",  "47.79", "2.82", "0.0161", 6)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    $rowIdx = $r + 2
    $q3.Cells.Item($rowIdx, 1).Value = $r
    $q3.Cells.Item($rowIdx, 2).Value = $row[0]
    $q3.Cells.Item($rowIdx, 3).Value = $row[1]
    $q3.Cells.Item($rowIdx, 4).Value = $row[2]
    $q3.Cells.Item($rowIdx, 5).Value = $row[3]
    $q3.Cells.Item($rowIdx, 6).Value = $row[4]
    $q3.Cells.Item($rowIdx, 7).Value = $row[5]
    $q3.Cells.Item($rowIdx, 8).Value = $row[6]
}

Write-Host "2022-Q3 sheet and summary row added"
